$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Summary sheet: update selection
# ------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("C5").Select()

# ------------------------------------------------------------------
# Repayment schedule sheet: data updates + remove column O cell
# definitions (O2:O8) while keeping column P untouched, then update
# selection.
# ------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Activate()

# Remove the O2:O8 cell records entirely (clear style back to default
# and blank the value so the cell is dropped on save).
$wsRepay.Range("O2:O8").Style = "Normal"
$wsRepay.Range("O2:O8").Value = ""

# Update row 5 values
$wsRepay.Range("I5").Value = 50
$wsRepay.Range("K5").Value = 937.72
$wsRepay.Range("P5").Value = 937.72

$wsRepay.Range("F22").Select()

# ------------------------------------------------------------------
# Transactions sheet: update values + selection
# ------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()

$wsTrans.Range("A2").Value = 16
$wsTrans.Range("A3").Value = 12

$wsTrans.Range("A2:L3").Select()
